# Adds a "2022-Q1" sheet (new fund-holdings snapshot) before "总计",
# and updates "总计" with a new leading row summarising 2022-Q1.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right before "总计"
# ------------------------------------------------------------------
$before = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($before)
$q1.Name = "2022-Q1"

# Base the header row + index-column formatting on "2021-Q4" (same layout/style)
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy($q1.Range("B1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund codes / percentages)
# -- format as Text first so leading zeros and exact decimals survive
$q1.Range("B2:B12").NumberFormat = "@"
$q1.Range("D2:G12").NumberFormat = "@"

# ------------------------------------------------------------------
# 2) Populate the 11 fund rows
# ------------------------------------------------------------------
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "360006"
$q1.Cells.Item(2, 3).Value = "光大保德信新增长混合"
$q1.Cells.Item(2, 4).Value = "21.71"
$q1.Cells.Item(2, 5).Value = "88.07"
$q1.Cells.Item(2, 6).Value = "3.44"
$q1.Cells.Item(2, 7).Value = "0.7468"
$q1.Cells.Item(2, 8).Value = 10
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "010330"
$q1.Cells.Item(3, 3).Value = "东吴兴享成长混合A"
$q1.Cells.Item(3, 4).Value = "11.63"
$q1.Cells.Item(3, 5).Value = "80.15"
$q1.Cells.Item(3, 6).Value = "5.21"
$q1.Cells.Item(3, 7).Value = "0.6059"
$q1.Cells.Item(3, 8).Value = 4
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = "011104"
$q1.Cells.Item(4, 3).Value = "光大保德信智能汽车主题股票"
$q1.Cells.Item(4, 4).Value = "10.77"
$q1.Cells.Item(4, 5).Value = "90.06"
$q1.Cells.Item(4, 6).Value = "4.33"
$q1.Cells.Item(4, 7).Value = "0.4663"
$q1.Cells.Item(4, 8).Value = 7
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = "001740"
$q1.Cells.Item(5, 3).Value = "光大保德信中国制造2025灵活配置混合"
$q1.Cells.Item(5, 4).Value = "11.43"
$q1.Cells.Item(5, 5).Value = "86.23"
$q1.Cells.Item(5, 6).Value = "3.26"
$q1.Cells.Item(5, 7).Value = "0.3726"
$q1.Cells.Item(5, 8).Value = 7
$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).Value = "010676"
$q1.Cells.Item(6, 3).Value = "光大保德信新机遇混合"
$q1.Cells.Item(6, 4).Value = "4.20"
$q1.Cells.Item(6, 5).Value = "85.57"
$q1.Cells.Item(6, 6).Value = "4.44"
$q1.Cells.Item(6, 7).Value = "0.1865"
$q1.Cells.Item(6, 8).Value = 9
$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).Value = "360005"
$q1.Cells.Item(7, 3).Value = "光大保德信红利混合"
$q1.Cells.Item(7, 4).Value = "4.34"
$q1.Cells.Item(7, 5).Value = "85.30"
$q1.Cells.Item(7, 6).Value = "2.84"
$q1.Cells.Item(7, 7).Value = "0.1233"
$q1.Cells.Item(7, 8).Value = 7
$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).Value = "006233"
$q1.Cells.Item(8, 3).Value = "万家汽车新趋势混合A"
$q1.Cells.Item(8, 4).Value = "4.33"
$q1.Cells.Item(8, 5).Value = "92.45"
$q1.Cells.Item(8, 6).Value = "2.56"
$q1.Cells.Item(8, 7).Value = "0.1108"
$q1.Cells.Item(8, 8).Value = 6
$q1.Cells.Item(9, 1).Value = 7
$q1.Cells.Item(9, 2).Value = "006234"
$q1.Cells.Item(9, 3).Value = "万家汽车新趋势混合C"
$q1.Cells.Item(9, 4).Value = "2.52"
$q1.Cells.Item(9, 5).Value = "92.45"
$q1.Cells.Item(9, 6).Value = "2.56"
$q1.Cells.Item(9, 7).Value = "0.0645"
$q1.Cells.Item(9, 8).Value = 6
$q1.Cells.Item(10, 1).Value = 8
$q1.Cells.Item(10, 2).Value = "011462"
$q1.Cells.Item(10, 3).Value = "东吴兴享成长混合C"
$q1.Cells.Item(10, 4).Value = "0.33"
$q1.Cells.Item(10, 5).Value = "80.15"
$q1.Cells.Item(10, 6).Value = "5.21"
$q1.Cells.Item(10, 7).Value = "0.0172"
$q1.Cells.Item(10, 8).Value = 4
$q1.Cells.Item(11, 1).Value = 9
$q1.Cells.Item(11, 2).Value = "673081"
$q1.Cells.Item(11, 3).Value = "西部利得祥运灵活配置混合A"
$q1.Cells.Item(11, 4).Value = "0.14"
$q1.Cells.Item(11, 5).Value = "84.19"
$q1.Cells.Item(11, 6).Value = "4.50"
$q1.Cells.Item(11, 7).Value = "0.0063"
$q1.Cells.Item(11, 8).Value = 3
$q1.Cells.Item(12, 1).Value = 10
$q1.Cells.Item(12, 2).Value = "673083"
$q1.Cells.Item(12, 3).Value = "西部利得祥运灵活配置混合C"
$q1.Cells.Item(12, 4).Value = "0.06"
$q1.Cells.Item(12, 5).Value = "84.19"
$q1.Cells.Item(12, 6).Value = "4.50"
$q1.Cells.Item(12, 7).Value = "0.0027"
$q1.Cells.Item(12, 8).Value = 3

# Apply the bold index-column style to A2:A12 (matches column A elsewhere in the workbook)
$src.Range("A2").Copy()
$q1.Range("A2:A12").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Update "总计": insert a new leading row for 2022-Q1
#    (re-fetch the sheet by name now that sheet positions have shifted)
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 2.7

# Re-apply the bold index-column style to the new A2 (Insert() does not carry it over)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows that shifted down, 1..3
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3

# ------------------------------------------------------------------
# 4) Restore the original active sheet/tab selection
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
